# Login TC - Get Expected Text from Excel
# Adds a new (blank) "Sheet1" worksheet at the end of the workbook, and adds
# the expected login result strings (used by the Login test case) to the
# "User Details" sheet. Also moves the active/selected sheet back to
# "User Details".

$wb = $excel.ActiveWorkbook

# Worker Details / TimeSheet Details etc. stay untouched; grab the last
# existing sheet so the new sheet is appended after it.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)

# The "User Details" sheet is the first sheet in the workbook.
$ws = $wb.Worksheets.Item(1)

# Expected messages for the Login test case, written column-by-column
# (B first, then A) so the shared-string table is built in the same order
# the original change used.
$ws.Range("B3").Value = "Welcome to Payroll Application"
$ws.Range("B4").Value = "Incorrect username or password."
$ws.Range("A3").Value = "validLogin"
$ws.Range("A4").Value = "invalidLogin"

# Widen column A slightly to fit the new values, and leave the selection on
# the cell the author was last working in.
$ws.Columns.Item(1).ColumnWidth = 10.8
$ws.Range("E10").Select() | Out-Null

# Finally make "User Details" the active sheet/tab again.
$ws.Activate() | Out-Null
